$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 9
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 10
$ws.Range("C9").Value = 9
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 8
$ws.Range("C14").Value = 5

$ws.Range("D2:D15").WrapText = $true

$ws.Range("D6").Value = "lien a certaines pages non-fonctionnels"
$ws.Range("D9").Value = "certaines pages brisent le formattage du site"
$ws.Range("D12").Value = "très bonne description du site web, facile à trouver la page principale"
$ws.Range("D13").Value = "certaines pages on des commentaires pertinant, tandis que d'autres n'ont rien"
$ws.Range("D14").Value = "*petit problème de la part de sidi, mais a été résolus"

$ws.Range("C6").Select()

$wb.Save()
